# Auto-generated Excel COM-interop script
# Applies the 2023-02-20 daily-data-refresh values (column J = 2023 YTD totals,
# plus a couple of very small historical H/I corrections) across all affected sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 859
$ws.Range("J3").Value = 947
$ws.Range("H4").Value = 1686
$ws.Range("I4").Value = 1751
$ws.Range("J4").Value = 203
$ws.Range("J5").Value = 70
$ws.Range("J6").Value = 1345
$ws.Range("H7").Value = 25998
$ws.Range("I7").Value = 26183
$ws.Range("J7").Value = 3424

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 39
$ws.Range("J2").Value = 13

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 31
$ws.Range("J3").Value = 38
$ws.Range("H4").Value = 32
$ws.Range("H7").Value = 941
$ws.Range("J7").Value = 117

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 39

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 28
$ws.Range("J7").Value = 97
$ws.Range("J8").Value = 220
$ws.Range("J15").Value = 38
$ws.Range("J18").Value = 56
$ws.Range("J19").Value = 111
$ws.Range("J20").Value = 75
$ws.Range("J27").Value = 16
$ws.Range("J29").Value = 176
$ws.Range("J33").Value = 146
$ws.Range("J36").Value = 52
$ws.Range("H37").Value = 941
$ws.Range("J37").Value = 117
$ws.Range("J42").Value = 153
$ws.Range("J47").Value = 31
$ws.Range("I48").Value = 328
$ws.Range("J48").Value = 21
$ws.Range("J51").Value = 43
$ws.Range("J52").Value = 76
$ws.Range("J55").Value = 42
$ws.Range("J60").Value = 19
$ws.Range("J64").Value = 21
$ws.Range("J65").Value = 91
$ws.Range("J74").Value = 5
$ws.Range("J79").Value = 107
$ws.Range("J83").Value = 78
$ws.Range("J84").Value = 39
$ws.Range("J85").Value = 139
$ws.Range("J89").Value = 39
$ws.Range("J91").Value = 45
$ws.Range("J95").Value = 64
$ws.Range("J96").Value = 38
$ws.Range("H101").Value = 25998
$ws.Range("I101").Value = 26183
$ws.Range("J101").Value = 3424

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J3").Value = 25
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 78
$ws.Range("J2").Value = 23

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 64

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 32
$ws.Range("J3").Value = 39
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 146

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 21
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 111

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J3").Value = 4

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I4").Value = 42
$ws.Range("I7").Value = 328
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 44
$ws.Range("J7").Value = 139

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 24
$ws.Range("J6").Value = 96
$ws.Range("J7").Value = 153

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J2").Value = 13
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 29
$ws.Range("J3").Value = 33
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 107

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J2").Value = 3
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 56

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J2").Value = 5
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 13
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 65
$ws.Range("J3").Value = 74
$ws.Range("J4").Value = 12
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 220

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J2").Value = 2
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J4").Value = 9
$ws.Range("J5").Value = 2
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 19

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 27
$ws.Range("J7").Value = 97

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("J6").Value = 2
$ws.Range("J7").Value = 5
